# Issue #65 added version number to PlaylistDTO, and hence saved file
#
# Adds new rows (63, 64, 65) to the "Issues" log sheet, marks row 62 as
# DONE (column D), scrolls/selects to reflect where the author left the
# sheet, and leaves the rest of the workbook untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")
$ws.Activate()

# Row 62 gets marked DONE (Status column D) - it was left blank before.
$ws.Range("D62").Value = "DONE"

# --- New row 63 --------------------------------------------------------
$ws.Range("A63").Value = 63
$ws.Range("E63").Value = "arch "
$ws.Range("F63").Value = "emulator to share same server as main"
$ws.Range("I63").Value = "means that there are not two ports"

# --- New row 64 --------------------------------------------------------
$ws.Range("A64").Value = 64
$ws.Range("F64").Value = "angular logging"

# --- New row 65 --------------------------------------------------------
$ws.Range("A65").Value = 65
$ws.Range("D65").Value = "DONE"
$ws.Range("F65").Value = "version numbers on saved playlist"

# Match the row heights of the other two-line wrapped rows (63 and 65 wrap
# onto a second line in the real workbook, row 64 stays single-line).
$ws.Range("A63:I63").RowHeight = 29
$ws.Range("A65:I65").RowHeight = 29

# Restore the filtered/frozen view the author ended up on: scrolled down
# so row 40 is the first visible row under the frozen header, with the
# cursor resting on E72.
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("E72").Select()
